$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.109.97"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.64%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.142.76"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.40"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.73%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.137.14"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.57%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.516"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.98%  "
$ws.Range("E10").Value = "  -2.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.23"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.35%  "
$ws.Range("E12").Value = "  -2.26%  "
$ws.Range("E13").Value = "  -3.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.19"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.659.73"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.42%  "
$ws.Range("E16").Value = "  +1.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.138.79"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.063.41"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.92%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.91%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "471.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.85%  "
$ws.Range("E21").Value = "  -2.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.698"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.66"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("E24").Value = "  -2.61%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.98"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.85%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("E27").Value = "  -1.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.95"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.12%  "
$ws.Range("E30").Value = "  +2.84%  "
$ws.Range("E31").Value = "  -0.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.76"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.107"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.53"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.93%  "
$ws.Range("E35").Value = "  -2.75%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.78"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.48%  "
$ws.Range("E37").Value = "  -0.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0X0691"
$ws.Range("D38").Characters(4,1).Text = [char]0x2083
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -8.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0386"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "418.71"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.73"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.18"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.900.20"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.24%  "
$ws.Range("E44").Value = "  -6.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.261"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.74%  "
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.44"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.66%  "
$ws.Range("E49").Value = "  -0.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.25"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "120.55"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.85%  "
